# Weekly driver report update for 2025-04-28
# - Bad Drivers table: Critical Minutes / Good Roaming % refreshed for the week
# - Good Drivers table: new driver version (21.40.1.3) added at the top,
#   and client-count samples refreshed for the existing driver versions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" summary table -------------------------------------------
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 98.40000000000001

# Totals row
$ws.Range("C4").Value = 6

# --- "Good Drivers" table ----------------------------------------------------
# Insert a new row for driver version 21.40.1.3, which now has the lowest
# vintage and is placed first in the list; this shifts the remaining driver
# rows down by one.
$ws.Rows("12:12").Insert()

# Copy the formatting from the row directly below (the row that used to be
# row 12 before the insert) so the new row matches the table's styling.
$ws.Range("A13:E13").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the newly inserted row with the new driver's data.
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B12").Value = 11128
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = ""

# Refresh total-sample counts for the existing driver versions (rows shifted
# down by the insert above).
$ws.Range("B13").Value = 486214
$ws.Range("B14").Value = 79953
$ws.Range("B15").Value = 35355
$ws.Range("B16").Value = 65425
$ws.Range("B17").Value = 117653
